$d = $word.ActiveDocument

# --- 1) Paragraph 1: merge split runs / remove proofErr spell-check marks ---
$old1 = "En un lugar de la Mancha, de cuyo nombre no quiero acordarme, no ha mucho tiempo que vivía un hidalgo de los de lanza en astillero, adarga antigua, rocín flaco y galgo corredor. Una olla de algo más vaca que carnero, salpicón las más noches, duelos y quebrantos los sábados, lentejas los viernes, algún palomino de añadidura los domingos, consumían las tres partes de su hacienda. El resto della concluían sayo de velarte, calzas de velludo para las fiestas, con sus pantuflos de lo mesmo, y los días de entresemana se honraba con su vellorí de lo más fino. Tenía en su casa una ama que pasaba de los cuarenta, y una sobrina que no llegaba a los veinte, y un mozo de campo y plaza, que así ensillaba el rocín como tomaba la podadera. Frisaba la edad de nuestro hidalgo con los cincuenta años; era de complexión recia, seco de carnes, enjuto de rostro, gran madrugador y amigo de la caza. Quieren decir que tenía el sobrenombre de Quijada, o Quesada, que en esto hay alguna diferencia en los autores que deste caso escriben; aunque por conjeturas verosímiles se deja entender que se llamaba Quijana. Pero esto importa poco a nuestro cuento: basta que en la narración del no se salga un punto de la verdad."

$p1 = $d.Paragraphs(1)
$p1.Range.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $old1, 2) | Out-Null

# --- 2) Paragraph 3: merge split runs / remove proofErr spell-check marks ---
$old2 = "Es, pues, de saber que este sobredicho hidalgo, los ratos que estaba ocioso, que eran los más del año, se daba a leer libros de caballerías, con tanta afición y gusto, que olvidó casi de todo punto el ejercicio de la caza, y aun la administración de su hacienda; y llegó a tanto su curiosidad y desatino en esto, que vendió muchas hanegas de tierra de sembradura para comprar libros de caballerías en que leer, y así, llevó a su casa todos cuantos pudo haber dellos; y de todos, ningunos le parecían tan bien como los que compuso el famoso Feliciano de Silva; porque la claridad de su prosa y aquellas entricadas razones suyas le parecían de perlas, y más cuando llegaba a leer aquellos requiebros y cartas de desafíos, donde en muchas partes hallaba escrito: «La razón de la sinrazón que a mi razón se hace, de tal manera mi razón enflaquece, que con razón me quejo de la vuestra fermosura». Y también cuando leía: «... los altos cielos que de vuestra divinidad divinamente con las estrellas os fortifican, y os hacen merecedora del merecimiento que merece la vuestra grandeza»."

$p2 = $d.Paragraphs(3)
$p2.Range.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $old2, 2) | Out-Null

# --- 3) Last (empty) paragraph gains the commit-message text, in Courier New ---
$pLast = $d.Paragraphs($d.Paragraphs.Count)
$markRange = $pLast.Range.Duplicate
$markRange.Font.NameAscii = "Courier New"
$markRange.Font.NameOther = "Courier New"
$markRange.Font.NameBi = "Courier New"

$pLast.Range.InsertBefore("Modificando quijote.docx.")

$pLast2 = $d.Paragraphs($d.Paragraphs.Count)
$runRange = $pLast2.Range.Duplicate
$runRange.MoveEnd(1, -1) | Out-Null
$runRange.Font.NameAscii = "Courier New"
$runRange.Font.NameOther = "Courier New"
$runRange.Font.NameBi = "Courier New"

Write-Output "ok"
